# Adds new translation rows (settings / account / notifications / privacy /
# general / help / selectDistrict screens) to the translations sheet.
# Columns: A=key, B=en, C=hi, D=te. Mirrors the commit "search screens are changed".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(403, 1).Value = 'settings_title'
$ws.Cells.Item(403, 2).Value = 'Settings'
$ws.Cells.Item(403, 3).Value = 'सेटिंग्स'
$ws.Cells.Item(403, 4).Value = 'సెట్టింగ్‌లు'
$ws.Cells.Item(404, 1).Value = 'manage_preferences'
$ws.Cells.Item(404, 2).Value = 'Manage your app preferences'
$ws.Cells.Item(404, 3).Value = 'अपने ऐप की प्राथमिकताएँ प्रबंधित करें'
$ws.Cells.Item(404, 4).Value = 'మీ యాప్ అభిరుచులను నిర్వహించండి'

$ws.Cells.Item(406, 1).Value = 'account_title'
$ws.Cells.Item(406, 2).Value = 'Account'
$ws.Cells.Item(406, 3).Value = 'खाता'
$ws.Cells.Item(406, 4).Value = 'ఖాతా'
$ws.Cells.Item(407, 1).Value = 'edit_profile'
$ws.Cells.Item(407, 2).Value = 'Edit Profile'
$ws.Cells.Item(407, 3).Value = 'प्रोफ़ाइल संपादित करें'
$ws.Cells.Item(407, 4).Value = 'ప్రొఫైల్ ఎడిట్ చేయండి'
$ws.Cells.Item(408, 1).Value = 'change_password'
$ws.Cells.Item(408, 2).Value = 'Change Password'
$ws.Cells.Item(408, 3).Value = 'पासवर्ड बदलें'
$ws.Cells.Item(408, 4).Value = 'పాస్‌వర్డ్ మార్చండి'
$ws.Cells.Item(409, 1).Value = 'verification_status'
$ws.Cells.Item(409, 2).Value = 'Verification Status'
$ws.Cells.Item(409, 3).Value = 'सत्यापन स्थिति'
$ws.Cells.Item(409, 4).Value = 'తనీకరణ స్థితి'
$ws.Cells.Item(410, 1).Value = 'verified'
$ws.Cells.Item(410, 2).Value = 'Verified'
$ws.Cells.Item(410, 3).Value = 'सत्यापित'
$ws.Cells.Item(410, 4).Value = 'తనఖీ చేయబడింది'

$ws.Cells.Item(412, 1).Value = 'notifications_title'
$ws.Cells.Item(412, 2).Value = 'Notifications'
$ws.Cells.Item(412, 3).Value = 'सूचनाएँ'
$ws.Cells.Item(412, 4).Value = 'నోటిఫికేషన్లు'
$ws.Cells.Item(413, 1).Value = 'property_alerts'
$ws.Cells.Item(413, 2).Value = 'Property Alerts'
$ws.Cells.Item(413, 3).Value = 'प्रॉपर्टी अलर्ट'
$ws.Cells.Item(413, 4).Value = 'ప్రాపర్టీ అలర్ట్స్'
$ws.Cells.Item(414, 1).Value = 'price_changes'
$ws.Cells.Item(414, 2).Value = 'Price Changes'
$ws.Cells.Item(414, 3).Value = 'कीमत में बदलाव'
$ws.Cells.Item(414, 4).Value = 'ధర మార్పులు'
$ws.Cells.Item(415, 1).Value = 'messages'
$ws.Cells.Item(415, 2).Value = 'Messages'
$ws.Cells.Item(415, 3).Value = 'संदेश'
$ws.Cells.Item(415, 4).Value = 'సందేశాలు'
$ws.Cells.Item(416, 1).Value = 'marketing_emails'
$ws.Cells.Item(416, 2).Value = 'Marketing Emails'
$ws.Cells.Item(416, 3).Value = 'मार्केटिंग ईमेल'
$ws.Cells.Item(416, 4).Value = 'మార్కెటింగ్ ఇమెయిల్స్'

$ws.Cells.Item(418, 1).Value = 'privacy_title'
$ws.Cells.Item(418, 2).Value = 'Privacy'
$ws.Cells.Item(418, 3).Value = 'गोपनीयता'
$ws.Cells.Item(418, 4).Value = 'గోప్యత'
$ws.Cells.Item(419, 1).Value = 'profile_visible'
$ws.Cells.Item(419, 2).Value = 'Profile Visible'
$ws.Cells.Item(419, 3).Value = 'प्रोफ़ाइल दिखाई दे'
$ws.Cells.Item(419, 4).Value = 'ప్రొఫైల్ కనిపించాలి'
$ws.Cells.Item(420, 1).Value = 'show_activity'
$ws.Cells.Item(420, 2).Value = 'Show Activity'
$ws.Cells.Item(420, 3).Value = 'गतिविधि दिखाएँ'
$ws.Cells.Item(420, 4).Value = 'కార్యకలాపం చూపించండి'
$ws.Cells.Item(421, 1).Value = 'allow_messages'
$ws.Cells.Item(421, 2).Value = 'Allow Messages'
$ws.Cells.Item(421, 3).Value = 'संदेशों की अनुमति दें'
$ws.Cells.Item(421, 4).Value = 'సందేశాలను అనుమతించండి'

$ws.Cells.Item(423, 1).Value = 'general_title'
$ws.Cells.Item(423, 2).Value = 'General'
$ws.Cells.Item(423, 3).Value = 'सामान्य'
$ws.Cells.Item(423, 4).Value = 'సాధారణ'
$ws.Cells.Item(424, 1).Value = 'language'
$ws.Cells.Item(424, 2).Value = 'Language'
$ws.Cells.Item(424, 3).Value = 'भाषा'
$ws.Cells.Item(424, 4).Value = 'భాష'
$ws.Cells.Item(425, 1).Value = 'currency'
$ws.Cells.Item(425, 2).Value = 'Currency'
$ws.Cells.Item(425, 3).Value = 'मुद्रा'
$ws.Cells.Item(425, 4).Value = 'కరెన్సీ'
$ws.Cells.Item(426, 1).Value = 'help_support'
$ws.Cells.Item(426, 2).Value = 'Help & Support'
$ws.Cells.Item(426, 3).Value = 'सहायता और समर्थन'
$ws.Cells.Item(426, 4).Value = 'సహాయం & మద్దతు'

$ws.Cells.Item(428, 1).Value = 'selectDistrict.title'
$ws.Cells.Item(428, 2).Value = 'Select District in Andhra'
$ws.Cells.Item(428, 3).Value = 'आंध्र में जिला चुनें'
$ws.Cells.Item(428, 4).Value = 'ఆంధ్రప్రదేశ్ జిల్లాను ఎంచుకోండి'
$ws.Cells.Item(429, 1).Value = 'selectDistrict.searchPlaceholder'
$ws.Cells.Item(429, 2).Value = 'Search districts in Andhra...'
$ws.Cells.Item(429, 3).Value = 'आंध्र के जिलों को खोजें...'
$ws.Cells.Item(429, 4).Value = 'ఆంధ్రా జిల్లాలను శోధించండి...'
$ws.Cells.Item(430, 1).Value = 'selectDistrict.noResults'
$ws.Cells.Item(430, 2).Value = 'No districts found for "{query}"'
$ws.Cells.Item(430, 3).Value = '{query} के लिए कोई जिला नहीं मिला'
$ws.Cells.Item(430, 4).Value = '{query} కి సంబంధించిన జిల్లాలు లభించలేదు'
$ws.Cells.Item(431, 1).Value = 'selectDistrict.propertiesAvailable'
$ws.Cells.Item(431, 2).Value = 'properties available'
$ws.Cells.Item(431, 3).Value = 'प्रॉपर्टीज उपलब्ध'
$ws.Cells.Item(431, 4).Value = 'ప్రాపర్టీలు అందుబాటులో ఉన్నాయి'

# Match the author's final selection (one blank row below the last new block).
[void]$ws.Range("A435").Select()
